$d = $word.ActiveDocument

# Locate the target bullet paragraph by its unique text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Ubicar un bot.n de guardar a la vista sin necesidad de ir hasta el men. archivo") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Include the paragraph mark itself so the highlight also lands on
    # the <w:pPr><w:rPr> (paragraph mark run properties), matching the
    # diff which highlights both the run and the paragraph mark.
    $r = $target.Range
    $r.HighlightColorIndex = 7   # wdYellow
}
